# Auto-generated edit script applying numeric updates per the commit diff
# Sheets map: 1=ALC 2=ARM 3=BSM 4=CRP 5=CUL 6=GSM 7=LTW 8=WVR
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)   # ALC
# Row 6
$ws.Cells.Item(6, 8).Value = 840.8333  # H6: 618.73334 -> 840.8333
$ws.Cells.Item(6, 9).Value = 899.0909  # I6: 618.73334 -> 899.0909
$ws.Cells.Item(6, 10).Value = 200  # J6: 0 -> 200
$ws.Cells.Item(6, 11).Value = 2697.2727  # K6: 1856.20002 -> 2697.2727
$ws.Cells.Item(6, 12).Value = 600  # L6: 0 -> 600
$ws.Cells.Item(6, 13).Value = -2585.2727  # M6: -1744.20002 -> -2585.2727
$ws.Cells.Item(6, 14).Value = -824  # N6: None -> -824

# Row 87
$ws.Cells.Item(87, 8).Value = 57995  # H87: 53995 -> 57995
$ws.Cells.Item(87, 10).Value = 57995  # J87: 53995 -> 57995
$ws.Cells.Item(87, 12).Value = 57995  # L87: 53995 -> 57995
$ws.Cells.Item(87, 14).Value = -60491  # N87: -56491 -> -60491

# Row 90
$ws.Cells.Item(90, 8).Value = 57995  # H90: 53995 -> 57995
$ws.Cells.Item(90, 10).Value = 57995  # J90: 53995 -> 57995
$ws.Cells.Item(90, 12).Value = 173985  # L90: 161985 -> 173985
$ws.Cells.Item(90, 14).Value = -186465  # N90: -174465 -> -186465

# Row 100
$ws.Cells.Item(100, 8).Value = 1955  # H100: 2204.4375 -> 1955
$ws.Cells.Item(100, 9).Value = 1345.5  # I100: 1492.5 -> 1345.5
$ws.Cells.Item(100, 10).Value = 3986.6667  # J100: 3391 -> 3986.6667
$ws.Cells.Item(100, 11).Value = 1345.5  # K100: 1492.5 -> 1345.5
$ws.Cells.Item(100, 12).Value = 3986.6667  # L100: 3391 -> 3986.6667
$ws.Cells.Item(100, 13).Value = -804.5  # M100: -951.5 -> -804.5
$ws.Cells.Item(100, 14).Value = -5068.6667  # N100: -4473 -> -5068.6667

# Row 111
$ws.Cells.Item(111, 8).Value = 652  # H111: 645 -> 652
$ws.Cells.Item(111, 9).Value = 570  # I111: 578 -> 570
$ws.Cells.Item(111, 10).Value = 843.3333  # J111: 980 -> 843.3333
$ws.Cells.Item(111, 11).Value = 1710  # K111: 1734 -> 1710
$ws.Cells.Item(111, 12).Value = 2529.9999  # L111: 2940 -> 2529.9999
$ws.Cells.Item(111, 13).Value = 1357  # M111: 1333 -> 1357
$ws.Cells.Item(111, 14).Value = -8663.999899999999  # N111: -9074 -> -8663.999899999999

# Row 113
$ws.Cells.Item(113, 8).Value = 2433.7827  # H113: 2518.1428 -> 2433.7827
$ws.Cells.Item(113, 9).Value = 1626.9  # I113: 1679.8889 -> 1626.9
$ws.Cells.Item(113, 10).Value = 3054.4614  # J113: 3146.8333 -> 3054.4614
$ws.Cells.Item(113, 11).Value = 1626.9  # K113: 1679.8889 -> 1626.9
$ws.Cells.Item(113, 12).Value = 3054.4614  # L113: 3146.8333 -> 3054.4614
$ws.Cells.Item(113, 13).Value = 1627.1  # M113: 1574.1111 -> 1627.1
$ws.Cells.Item(113, 14).Value = -9562.4614  # N113: -9654.8333 -> -9562.4614

# Row 118
$ws.Cells.Item(118, 8).Value = 862.2826  # H118: 1245.2128 -> 862.2826
$ws.Cells.Item(118, 9).Value = 419.70587  # I118: 446.07144 -> 419.70587
$ws.Cells.Item(118, 10).Value = 1121.7241  # J118: 1584.2424 -> 1121.7241
$ws.Cells.Item(118, 11).Value = 1259.11761  # K118: 1338.21432 -> 1259.11761
$ws.Cells.Item(118, 12).Value = 3365.1723  # L118: 4752.7272 -> 3365.1723
$ws.Cells.Item(118, 13).Value = 397.88239  # M118: 318.78568 -> 397.88239
$ws.Cells.Item(118, 14).Value = -6679.1723  # N118: -8066.7272 -> -6679.1723

# Row 137
$ws.Cells.Item(137, 8).Value = 3285.3572  # H137: 3504.3208 -> 3285.3572
$ws.Cells.Item(137, 9).Value = 1102.4814  # I137: 1191.6364 -> 1102.4814
$ws.Cells.Item(137, 10).Value = 5317.6895  # J137: 5145.5806 -> 5317.6895
$ws.Cells.Item(137, 11).Value = 3307.4442  # K137: 3574.9092 -> 3307.4442
$ws.Cells.Item(137, 12).Value = 15953.0685  # L137: 15436.7418 -> 15953.0685
$ws.Cells.Item(137, 13).Value = -757.4441999999999  # M137: -1024.9092 -> -757.4441999999999
$ws.Cells.Item(137, 14).Value = -21053.0685  # N137: -20536.7418 -> -21053.0685


# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)   # ARM
# Row 32
$ws.Cells.Item(32, 8).Value = 3123.12  # H32: 4764.227 -> 3123.12
$ws.Cells.Item(32, 9).Value = 2201.25  # I32: 3909.5325 -> 2201.25
$ws.Cells.Item(32, 10).Value = 9883.5  # J32: 10747.091 -> 9883.5
$ws.Cells.Item(32, 11).Value = 2201.25  # K32: 3909.5325 -> 2201.25
$ws.Cells.Item(32, 12).Value = 9883.5  # L32: 10747.091 -> 9883.5
$ws.Cells.Item(32, 13).Value = -1914.25  # M32: -3622.5325 -> -1914.25
$ws.Cells.Item(32, 14).Value = -10457.5  # N32: -11321.091 -> -10457.5

# Row 82
$ws.Cells.Item(82, 8).Value = 38436.2  # H82: 43333.332 -> 38436.2
$ws.Cells.Item(82, 10).Value = 38436.2  # J82: 43333.332 -> 38436.2
$ws.Cells.Item(82, 12).Value = 38436.2  # L82: 43333.332 -> 38436.2
$ws.Cells.Item(82, 14).Value = -39158.2  # N82: -44055.332 -> -39158.2

# Row 85
$ws.Cells.Item(85, 8).Value = 38436.2  # H85: 43333.332 -> 38436.2
$ws.Cells.Item(85, 10).Value = 38436.2  # J85: 43333.332 -> 38436.2
$ws.Cells.Item(85, 12).Value = 38436.2  # L85: 43333.332 -> 38436.2
$ws.Cells.Item(85, 14).Value = -40932.2  # N85: -45829.332 -> -40932.2


# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)   # BSM
# Row 38
$ws.Cells.Item(38, 8).Value = 8036  # H38: 0 -> 8036
$ws.Cells.Item(38, 10).Value = 8036  # J38: 0 -> 8036
$ws.Cells.Item(38, 12).Value = 8036  # L38: 0 -> 8036
$ws.Cells.Item(38, 14).Value = -8868  # N38: None -> -8868


# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)   # CRP
# Row 31
$ws.Cells.Item(31, 8).Value = 34512.26  # H31: 42932.33 -> 34512.26
$ws.Cells.Item(31, 9).Value = 1000.675  # I31: 1148.1852 -> 1000.675
$ws.Cells.Item(31, 10).Value = 88130.8  # J31: 88059.2 -> 88130.8
$ws.Cells.Item(31, 11).Value = 1000.675  # K31: 1148.1852 -> 1000.675
$ws.Cells.Item(31, 12).Value = 88130.8  # L31: 88059.2 -> 88130.8
$ws.Cells.Item(31, 13).Value = -705.675  # M31: -853.1851999999999 -> -705.675
$ws.Cells.Item(31, 14).Value = -88720.8  # N31: -88649.2 -> -88720.8

# Row 34
$ws.Cells.Item(34, 8).Value = 34512.26  # H34: 42932.33 -> 34512.26
$ws.Cells.Item(34, 9).Value = 1000.675  # I34: 1148.1852 -> 1000.675
$ws.Cells.Item(34, 10).Value = 88130.8  # J34: 88059.2 -> 88130.8
$ws.Cells.Item(34, 11).Value = 1000.675  # K34: 1148.1852 -> 1000.675
$ws.Cells.Item(34, 12).Value = 88130.8  # L34: 88059.2 -> 88130.8
$ws.Cells.Item(34, 13).Value = -798.675  # M34: -946.1851999999999 -> -798.675
$ws.Cells.Item(34, 14).Value = -88534.8  # N34: -88463.2 -> -88534.8

# Row 109
$ws.Cells.Item(109, 8).Value = 23000  # H109: 19000 -> 23000
$ws.Cells.Item(109, 10).Value = 23000  # J109: 19000 -> 23000
$ws.Cells.Item(109, 12).Value = 23000  # L109: 19000 -> 23000
$ws.Cells.Item(109, 14).Value = -25080  # N109: -21080 -> -25080

# Row 132
$ws.Cells.Item(132, 8).Value = 14087361  # H132: 16952446 -> 14087361
$ws.Cells.Item(132, 9).Value = 19234148  # I132: 23259762 -> 19234148
$ws.Cells.Item(132, 10).Value = 1418.1052  # J132: 1533.125 -> 1418.1052
$ws.Cells.Item(132, 11).Value = 57702444  # K132: 69779286 -> 57702444
$ws.Cells.Item(132, 12).Value = 4254.3156  # L132: 4599.375 -> 4254.3156
$ws.Cells.Item(132, 13).Value = -57699914  # M132: -69776756 -> -57699914
$ws.Cells.Item(132, 14).Value = -9314.3156  # N132: -9659.375 -> -9314.3156


# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)   # CUL
# Row 7
$ws.Cells.Item(7, 8).Value = 721.913  # H7: 799.3 -> 721.913
$ws.Cells.Item(7, 9).Value = 20  # I7: 366.66666 -> 20
$ws.Cells.Item(7, 10).Value = 788.7619  # J7: 875.64703 -> 788.7619
$ws.Cells.Item(7, 11).Value = 60  # K7: 1099.99998 -> 60
$ws.Cells.Item(7, 12).Value = 2366.2857  # L7: 2626.94109 -> 2366.2857
$ws.Cells.Item(7, 13).Value = 52  # M7: -987.9999800000001 -> 52
$ws.Cells.Item(7, 14).Value = -2590.2857  # N7: -2850.94109 -> -2590.2857

# Row 17
$ws.Cells.Item(17, 8).Value = 8400.666999999999  # H17: 6655.091 -> 8400.666999999999
$ws.Cells.Item(17, 9).Value = 500  # I17: 333.33334 -> 500
$ws.Cells.Item(17, 10).Value = 9388.25  # J17: 9025.75 -> 9388.25
$ws.Cells.Item(17, 11).Value = 1500  # K17: 1000.00002 -> 1500
$ws.Cells.Item(17, 12).Value = 28164.75  # L17: 27077.25 -> 28164.75
$ws.Cells.Item(17, 13).Value = -1331  # M17: -831.0000200000001 -> -1331
$ws.Cells.Item(17, 14).Value = -28502.75  # N17: -27415.25 -> -28502.75

# Row 34
$ws.Cells.Item(34, 8).Value = 220.5  # H34: 620.38464 -> 220.5
$ws.Cells.Item(34, 9).Value = 220.5  # I34: 192.4 -> 220.5
$ws.Cells.Item(34, 10).Value = 0  # J34: 887.875 -> 0
$ws.Cells.Item(34, 11).Value = 661.5  # K34: 577.2 -> 661.5
$ws.Cells.Item(34, 12).Value = 0  # L34: 2663.625 -> 0
$ws.Cells.Item(34, 13).Value = -577.5  # M34: -493.2 -> -577.5
$ws.Cells.Item(34, 14).ClearContents()  # N34: was -2831.625 -> removed

# Row 39
$ws.Cells.Item(39, 8).Value = 3487.5  # H39: 3536.3635 -> 3487.5
$ws.Cells.Item(39, 10).Value = 3857.1428  # J39: 3800 -> 3857.1428
$ws.Cells.Item(39, 12).Value = 11571.4284  # L39: 11400 -> 11571.4284
$ws.Cells.Item(39, 14).Value = -12159.4284  # N39: -11988 -> -12159.4284

# Row 55
$ws.Cells.Item(55, 8).Value = 2584.2856  # H55: 2820 -> 2584.2856
$ws.Cells.Item(55, 9).Value = 550  # I55: 850 -> 550
$ws.Cells.Item(55, 10).Value = 3398  # J55: 3312.5 -> 3398
$ws.Cells.Item(55, 11).Value = 1650  # K55: 2550 -> 1650
$ws.Cells.Item(55, 12).Value = 10194  # L55: 9937.5 -> 10194
$ws.Cells.Item(55, 13).Value = -1473  # M55: -2373 -> -1473
$ws.Cells.Item(55, 14).Value = -10548  # N55: -10291.5 -> -10548

# Row 80
$ws.Cells.Item(80, 8).Value = 3230.8462  # H80: 49113.684 -> 3230.8462
$ws.Cells.Item(80, 9).Value = 2771.7144  # I80: 202900.2 -> 2771.7144
$ws.Cells.Item(80, 10).Value = 3400  # J80: 3882.353 -> 3400
$ws.Cells.Item(80, 11).Value = 8315.143199999999  # K80: 608700.6000000001 -> 8315.143199999999
$ws.Cells.Item(80, 12).Value = 10200  # L80: 11647.059 -> 10200
$ws.Cells.Item(80, 13).Value = -7379.143199999999  # M80: -607764.6000000001 -> -7379.143199999999
$ws.Cells.Item(80, 14).Value = -12072  # N80: -13519.059 -> -12072

# Row 83
$ws.Cells.Item(83, 8).Value = 3230.8462  # H83: 49113.684 -> 3230.8462
$ws.Cells.Item(83, 9).Value = 2771.7144  # I83: 202900.2 -> 2771.7144
$ws.Cells.Item(83, 10).Value = 3400  # J83: 3882.353 -> 3400
$ws.Cells.Item(83, 11).Value = 24945.4296  # K83: 1826101.8 -> 24945.4296
$ws.Cells.Item(83, 12).Value = 30600  # L83: 34941.177 -> 30600
$ws.Cells.Item(83, 13).Value = -20265.4296  # M83: -1821421.8 -> -20265.4296
$ws.Cells.Item(83, 14).Value = -39960  # N83: -44301.177 -> -39960

# Row 92
$ws.Cells.Item(92, 8).Value = 352.5  # H92: 353.33334 -> 352.5
$ws.Cells.Item(92, 10).Value = 370  # J92: 380 -> 370
$ws.Cells.Item(92, 12).Value = 1110  # L92: 1140 -> 1110
$ws.Cells.Item(92, 14).Value = -3606  # N92: -3636 -> -3606


# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)   # GSM
# Row 80
$ws.Cells.Item(80, 8).Value = 2911.8125  # H80: 2596.7646 -> 2911.8125
$ws.Cells.Item(80, 9).Value = 2202.5  # I80: 2001.875 -> 2202.5
$ws.Cells.Item(80, 10).Value = 3013.1428  # J80: 3125.5557 -> 3013.1428
$ws.Cells.Item(80, 11).Value = 2202.5  # K80: 2001.875 -> 2202.5
$ws.Cells.Item(80, 12).Value = 3013.1428  # L80: 3125.5557 -> 3013.1428
$ws.Cells.Item(80, 13).Value = -1204.5  # M80: -1003.875 -> -1204.5
$ws.Cells.Item(80, 14).Value = -5009.1428  # N80: -5121.5557 -> -5009.1428

# Row 83
$ws.Cells.Item(83, 8).Value = 2911.8125  # H83: 2596.7646 -> 2911.8125
$ws.Cells.Item(83, 9).Value = 2202.5  # I83: 2001.875 -> 2202.5
$ws.Cells.Item(83, 10).Value = 3013.1428  # J83: 3125.5557 -> 3013.1428
$ws.Cells.Item(83, 11).Value = 11012.5  # K83: 10009.375 -> 11012.5
$ws.Cells.Item(83, 12).Value = 15065.714  # L83: 15627.7785 -> 15065.714
$ws.Cells.Item(83, 13).Value = -6020.5  # M83: -5017.375 -> -6020.5
$ws.Cells.Item(83, 14).Value = -25049.714  # N83: -25611.7785 -> -25049.714

# Row 132
$ws.Cells.Item(132, 8).Value = 2814.9055  # H132: 2858.1807 -> 2814.9055
$ws.Cells.Item(132, 9).Value = 3464.3  # I132: 3395.9019 -> 3464.3
$ws.Cells.Item(132, 10).Value = 1462  # J132: 1552.2858 -> 1462
$ws.Cells.Item(132, 11).Value = 10392.9  # K132: 10187.7057 -> 10392.9
$ws.Cells.Item(132, 12).Value = 4386  # L132: 4656.857400000001 -> 4386
$ws.Cells.Item(132, 13).Value = -7862.900000000001  # M132: -7657.705699999999 -> -7862.900000000001
$ws.Cells.Item(132, 14).Value = -9446  # N132: -9716.857400000001 -> -9446


# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)   # LTW
# Row 61
$ws.Cells.Item(61, 8).Value = 1538.1724  # H61: 1625.2693 -> 1538.1724
$ws.Cells.Item(61, 9).Value = 1264.5834  # I61: 1362.5 -> 1264.5834
$ws.Cells.Item(61, 10).Value = 1731.2941  # J61: 1789.5 -> 1731.2941
$ws.Cells.Item(61, 11).Value = 1264.5834  # K61: 1362.5 -> 1264.5834
$ws.Cells.Item(61, 12).Value = 1731.2941  # L61: 1789.5 -> 1731.2941
$ws.Cells.Item(61, 13).Value = -1062.5834  # M61: -1160.5 -> -1062.5834
$ws.Cells.Item(61, 14).Value = -2135.2941  # N61: -2193.5 -> -2135.2941

# Row 68
$ws.Cells.Item(68, 8).Value = 2048.5715  # H68: 1785.0333 -> 2048.5715
$ws.Cells.Item(68, 9).Value = 1965  # I68: 1788.1428 -> 1965
$ws.Cells.Item(68, 10).Value = 2160  # J68: 1777.7778 -> 2160
$ws.Cells.Item(68, 11).Value = 1965  # K68: 1788.1428 -> 1965
$ws.Cells.Item(68, 12).Value = 2160  # L68: 1777.7778 -> 2160
$ws.Cells.Item(68, 13).Value = -1216  # M68: -1039.1428 -> -1216
$ws.Cells.Item(68, 14).Value = -3658  # N68: -3275.7778 -> -3658

# Row 71
$ws.Cells.Item(71, 8).Value = 2048.5715  # H71: 1785.0333 -> 2048.5715
$ws.Cells.Item(71, 9).Value = 1965  # I71: 1788.1428 -> 1965
$ws.Cells.Item(71, 10).Value = 2160  # J71: 1777.7778 -> 2160
$ws.Cells.Item(71, 11).Value = 9825  # K71: 8940.714 -> 9825
$ws.Cells.Item(71, 12).Value = 10800  # L71: 8888.889000000001 -> 10800
$ws.Cells.Item(71, 13).Value = -6081  # M71: -5196.714 -> -6081
$ws.Cells.Item(71, 14).Value = -18288  # N71: -16376.889 -> -18288

# Row 82
$ws.Cells.Item(82, 8).Value = 1030  # H82: 936.8889 -> 1030
$ws.Cells.Item(82, 9).Value = 836  # I82: 850.6667 -> 836
$ws.Cells.Item(82, 10).Value = 2000  # J82: 980 -> 2000
$ws.Cells.Item(82, 11).Value = 836  # K82: 850.6667 -> 836
$ws.Cells.Item(82, 12).Value = 2000  # L82: 980 -> 2000
$ws.Cells.Item(82, 13).Value = -475  # M82: -489.6667 -> -475
$ws.Cells.Item(82, 14).Value = -2722  # N82: -1702 -> -2722

# Row 85
$ws.Cells.Item(85, 8).Value = 1030  # H85: 936.8889 -> 1030
$ws.Cells.Item(85, 9).Value = 836  # I85: 850.6667 -> 836
$ws.Cells.Item(85, 10).Value = 2000  # J85: 980 -> 2000
$ws.Cells.Item(85, 11).Value = 836  # K85: 850.6667 -> 836
$ws.Cells.Item(85, 12).Value = 2000  # L85: 980 -> 2000
$ws.Cells.Item(85, 13).Value = 412  # M85: 397.3333 -> 412
$ws.Cells.Item(85, 14).Value = -4496  # N85: -3476 -> -4496

# Row 113
$ws.Cells.Item(113, 8).Value = 1538.1724  # H113: 1625.2693 -> 1538.1724
$ws.Cells.Item(113, 9).Value = 1264.5834  # I113: 1362.5 -> 1264.5834
$ws.Cells.Item(113, 10).Value = 1731.2941  # J113: 1789.5 -> 1731.2941
$ws.Cells.Item(113, 11).Value = 1264.5834  # K113: 1362.5 -> 1264.5834
$ws.Cells.Item(113, 12).Value = 1731.2941  # L113: 1789.5 -> 1731.2941
$ws.Cells.Item(113, 13).Value = 905.4166  # M113: 807.5 -> 905.4166
$ws.Cells.Item(113, 14).Value = -6071.2941  # N113: -6129.5 -> -6071.2941

# Row 136
$ws.Cells.Item(136, 8).Value = 3848.825  # H136: 3771.878 -> 3848.825
$ws.Cells.Item(136, 9).Value = 1392.5  # I136: 1357.3871 -> 1392.5
$ws.Cells.Item(136, 10).Value = 11217.8  # J136: 11256.8 -> 11217.8
$ws.Cells.Item(136, 11).Value = 4177.5  # K136: 4072.1613 -> 4177.5
$ws.Cells.Item(136, 12).Value = 33653.39999999999  # L136: 33770.39999999999 -> 33653.39999999999
$ws.Cells.Item(136, 13).Value = -1627.5  # M136: -1522.1613 -> -1627.5
$ws.Cells.Item(136, 14).Value = -38753.39999999999  # N136: -38870.39999999999 -> -38753.39999999999


# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)   # WVR
# Row 62
$ws.Cells.Item(62, 8).Value = 2996.6667  # H62: 2983.3333 -> 2996.6667
$ws.Cells.Item(62, 9).Value = 2995  # I62: 3000 -> 2995
$ws.Cells.Item(62, 10).Value = 3000  # J62: 2966.6667 -> 3000
$ws.Cells.Item(62, 11).Value = 2995  # K62: 3000 -> 2995
$ws.Cells.Item(62, 12).Value = 3000  # L62: 2966.6667 -> 3000
$ws.Cells.Item(62, 13).Value = -2371  # M62: -2376 -> -2371
$ws.Cells.Item(62, 14).Value = -4248  # N62: -4214.6667 -> -4248

# Row 65
$ws.Cells.Item(65, 8).Value = 2996.6667  # H65: 2983.3333 -> 2996.6667
$ws.Cells.Item(65, 9).Value = 2995  # I65: 3000 -> 2995
$ws.Cells.Item(65, 10).Value = 3000  # J65: 2966.6667 -> 3000
$ws.Cells.Item(65, 11).Value = 14975  # K65: 15000 -> 14975
$ws.Cells.Item(65, 12).Value = 15000  # L65: 14833.3335 -> 15000
$ws.Cells.Item(65, 13).Value = -11855  # M65: -11880 -> -11855
$ws.Cells.Item(65, 14).Value = -21240  # N65: -21073.3335 -> -21240

# Row 107
$ws.Cells.Item(107, 8).Value = 259.35294  # H107: 285.7 -> 259.35294
$ws.Cells.Item(107, 9).Value = 251.35715  # I107: 275.875 -> 251.35715
$ws.Cells.Item(107, 10).Value = 296.66666  # J107: 325 -> 296.66666
$ws.Cells.Item(107, 11).Value = 754.0714499999999  # K107: 827.625 -> 754.0714499999999
$ws.Cells.Item(107, 12).Value = 889.9999799999999  # L107: 975 -> 889.9999799999999
$ws.Cells.Item(107, 13).Value = 1165.92855  # M107: 1092.375 -> 1165.92855
$ws.Cells.Item(107, 14).Value = -4729.99998  # N107: -4815 -> -4729.99998

# Row 132
$ws.Cells.Item(132, 8).Value = 766.59155  # H132: 3045.1785 -> 766.59155
$ws.Cells.Item(132, 9).Value = 583.283  # I132: 3912.5945 -> 583.283
$ws.Cells.Item(132, 10).Value = 1306.3334  # J132: 1356 -> 1306.3334
$ws.Cells.Item(132, 11).Value = 1749.849  # K132: 11737.7835 -> 1749.849
$ws.Cells.Item(132, 12).Value = 3919.0002  # L132: 4068 -> 3919.0002
$ws.Cells.Item(132, 13).Value = 780.1509999999998  # M132: -9207.783500000001 -> 780.1509999999998
$ws.Cells.Item(132, 14).Value = -8979.0002  # N132: -9128 -> -8979.0002

